$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price records in rows 2, 3 and 5 are re-ordered chronologically
# (row 4 stays as-is). This is effectively a 3-way rotation of the row data:
#   new row2 = old row3
#   new row3 = old row5
#   new row5 = old row2

# Capture the original values before overwriting anything.
$row2 = @{
    D = $ws.Range("D2").Value2
    I = $ws.Range("I2").Value2
    J = $ws.Range("J2").Value2
    K = $ws.Range("K2").Value2
    L = $ws.Range("L2").Value2
    M = $ws.Range("M2").Value2
    N = $ws.Range("N2").Value2
    P = $ws.Range("P2").Value2
    Q = $ws.Range("Q2").Value2
}

$row3 = @{
    D = $ws.Range("D3").Value2
    I = $ws.Range("I3").Value2
    J = $ws.Range("J3").Value2
    K = $ws.Range("K3").Value2
    L = $ws.Range("L3").Value2
    M = $ws.Range("M3").Value2
    N = $ws.Range("N3").Value2
    P = $ws.Range("P3").Value2
    Q = $ws.Range("Q3").Value2
}

$row5 = @{
    D = $ws.Range("D5").Value2
    I = $ws.Range("I5").Value2
    J = $ws.Range("J5").Value2
    K = $ws.Range("K5").Value2
    L = $ws.Range("L5").Value2
    M = $ws.Range("M5").Value2
    N = $ws.Range("N5").Value2
    P = $ws.Range("P5").Value2
    Q = $ws.Range("Q5").Value2
}

# Write old row3 data into row2
$ws.Range("D2").Value2 = $row3.D
$ws.Range("I2").Value2 = $row3.I
$ws.Range("J2").Value2 = $row3.J
$ws.Range("K2").Value2 = $row3.K
$ws.Range("L2").Value2 = $row3.L
$ws.Range("M2").Value2 = $row3.M
$ws.Range("N2").Value2 = $row3.N
$ws.Range("P2").Value2 = $row3.P
$ws.Range("Q2").Value2 = $row3.Q

# Write old row5 data into row3
$ws.Range("D3").Value2 = $row5.D
$ws.Range("I3").Value2 = $row5.I
$ws.Range("J3").Value2 = $row5.J
$ws.Range("K3").Value2 = $row5.K
$ws.Range("L3").Value2 = $row5.L
$ws.Range("M3").Value2 = $row5.M
$ws.Range("N3").Value2 = $row5.N
$ws.Range("P3").Value2 = $row5.P
$ws.Range("Q3").Value2 = $row5.Q

# Write old row2 data into row5
$ws.Range("D5").Value2 = $row2.D
$ws.Range("I5").Value2 = $row2.I
$ws.Range("J5").Value2 = $row2.J
$ws.Range("K5").Value2 = $row2.K
$ws.Range("L5").Value2 = $row2.L
$ws.Range("M5").Value2 = $row2.M
$ws.Range("N5").Value2 = $row2.N
$ws.Range("P5").Value2 = $row2.P
$ws.Range("Q5").Value2 = $row2.Q
